# Updates cryptos list values (Price and Volume(1h) columns) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.295.25"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.494.04"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'588.98"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'133.66"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'7.63"
$ws.Range("E9").Value = "  +6.37%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("D12").Value = "4.093.18"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "3.495.18"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "64.240.02"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'25.30"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "'5.77"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").Value = "'13.55"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "'385.89"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "3.633.64"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'74.16"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "'1.49"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").Value = "'8.15"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "'0.155"
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("D34").Value = "3.524.80"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'23.31"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").Value = "'5.32"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "'1.53"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "'165.74"
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("D41").Value = "'0.0785"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "'0.806"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'24.55"
$ws.Range("E45").Value = "  -3.62%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "2.434.95"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "'6.82"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").Value = "'0.914"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("E51").Value = "  -0.54%  "
